# The 'default' column is being removed from the survey / section1 /
# section2 sheets. In the original layout that is column R (18th column):
# column R held "default" and column S (19th) held "hideInContents".
# Deleting column R shifts "hideInContents" (and the boolean data below
# it) left into column R, which is exactly what the diff shows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("survey", "section1", "section2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(18).Delete() | Out-Null
    # Leave the new last column (old "hideInContents" column, now R)
    # selected, matching the saved selection state in the edited file.
    $ws.Columns.Item(18).Select() | Out-Null
}

# Restore the original active sheet/tab: the edited workbook has the
# "settings" sheet active (last tab) instead of "survey" (first tab).
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
